$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 509.9
$ws.Range("I2").Value = 309.8
$ws.Range("J2").Value = 710
$ws.Range("K2").Value = 309.8
$ws.Range("L2").Value = 710
$ws.Range("M2").Value = -196.8
$ws.Range("N2").Value = -936

# Row 33
$ws.Range("H33").Value = 127.6
$ws.Range("I33").Value = 102.42857
$ws.Range("J33").Value = 480
$ws.Range("K33").Value = 102.42857
$ws.Range("L33").Value = 480
$ws.Range("M33").Value = 126.57143
$ws.Range("N33").Value = -938

# Row 43
$ws.Range("H43").Value = 13914164
$ws.Range("I43").Value = 50300.5
$ws.Range("K43").Value = 50300.5
$ws.Range("M43").Value = -50231.5

# Row 61
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

# Row 86
$ws.Range("H86").Value = 7477.6665
$ws.Range("J86").Value = 9900
$ws.Range("L86").Value = 9900
$ws.Range("N86").Value = -12146

# Row 89
$ws.Range("H89").Value = 7477.6665
$ws.Range("J89").Value = 9900
$ws.Range("L89").Value = 49500
$ws.Range("N89").Value = -60732

# Row 116
$ws.Range("H116").Value = 3374.818
$ws.Range("I116").Value = 3160.3333
$ws.Range("J116").Value = 3750.1667
$ws.Range("K116").Value = 3160.3333
$ws.Range("L116").Value = 3750.1667
$ws.Range("M116").Value = 281.6667000000002
$ws.Range("N116").Value = -10634.1667

# Row 129
$ws.Range("H129").Value = 867.4946
$ws.Range("I129").Value = 384.54544
$ws.Range("J129").Value = 932.2805
$ws.Range("K129").Value = 1153.63632
$ws.Range("L129").Value = 2796.8415
$ws.Range("M129").Value = 3846.36368
$ws.Range("N129").Value = -12796.8415

# Row 132
$ws.Range("H132").Value = 12829972
$ws.Range("I132").Value = 15159000
$ws.Range("K132").Value = 45477000
$ws.Range("M132").Value = -45474470

# Row 138
$ws.Range("H138").Value = 1299.46
$ws.Range("I138").Value = 843.4706
$ws.Range("J138").Value = 1534.3636
$ws.Range("K138").Value = 2530.4118
$ws.Range("L138").Value = 4603.0908
$ws.Range("M138").Value = 2609.5882
$ws.Range("N138").Value = -14883.0908

# Row 141
$ws.Range("H141").Value = 932
$ws.Range("I141").Value = 932
$ws.Range("K141").Value = 2796
$ws.Range("M141").Value = 2384

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2962.353
$ws.Range("I32").Value = 3070.302
$ws.Range("K32").Value = 3070.302
$ws.Range("M32").Value = -2783.302

# Row 110
$ws.Range("H110").Value = 1241.85
$ws.Range("I110").Value = 747.4
$ws.Range("K110").Value = 747.4
$ws.Range("M110").Value = 1297.6

# Row 132
$ws.Range("H132").Value = 2141.077
$ws.Range("I132").Value = 1750.2858
$ws.Range("K132").Value = 5250.857400000001
$ws.Range("M132").Value = -2720.857400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 100991330
$ws.Range("I105").Value = 126238490
$ws.Range("K105").Value = 126238490
$ws.Range("M105").Value = -126236743

# Row 107
$ws.Range("H107").Value = 1658.0667
$ws.Range("I107").Value = 1390.4546
$ws.Range("J107").Value = 2394
$ws.Range("K107").Value = 1390.4546
$ws.Range("L107").Value = 2394
$ws.Range("M107").Value = 529.5454
$ws.Range("N107").Value = -6234

# Row 134
$ws.Range("H134").Value = 12877.167
$ws.Range("I134").Value = 1450
$ws.Range("K134").Value = 4350
$ws.Range("M134").Value = -1815

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1275.909
$ws.Range("I31").Value = 1275.909
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1275.909
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -980.9090000000001
$ws.Range("N31").ClearContents()

# Row 34
$ws.Range("H34").Value = 1275.909
$ws.Range("I34").Value = 1275.909
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1275.909
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1073.909
$ws.Range("N34").ClearContents()

# Row 39
$ws.Range("H39").Value = 724.75
$ws.Range("I39").Value = 724.75
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 724.75
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -333.75
$ws.Range("N39").ClearContents()

# Row 49
$ws.Range("H49").Value = 724.75
$ws.Range("I49").Value = 724.75
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 724.75
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -542.75
$ws.Range("N49").ClearContents()

# Row 132
$ws.Range("H132").Value = 4116.5366
$ws.Range("I132").Value = 4488.552
$ws.Range("J132").Value = 3217.5
$ws.Range("K132").Value = 13465.656
$ws.Range("L132").Value = 9652.5
$ws.Range("M132").Value = -10935.656
$ws.Range("N132").Value = -14712.5

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 587.2143
$ws.Range("J121").Value = 884.6
$ws.Range("L121").Value = 2653.8
$ws.Range("N121").Value = -5273.8

# Row 131
$ws.Range("H131").Value = 20001286
$ws.Range("I131").Value = 125000500
$ws.Range("J131").Value = 1435.5952
$ws.Range("K131").Value = 375001500
$ws.Range("L131").Value = 4306.7856
$ws.Range("M131").Value = -374996460
$ws.Range("N131").Value = -14386.7856

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 56256750
$ws.Range("I70").Value = 83338664
$ws.Range("K70").Value = 83338664
$ws.Range("M70").Value = -83338394

# Row 73
$ws.Range("H73").Value = 56256750
$ws.Range("I73").Value = 83338664
$ws.Range("K73").Value = 83338664
$ws.Range("M73").Value = -83337728

# Row 102
$ws.Range("H102").Value = 1973
$ws.Range("I102").Value = 2052.889
$ws.Range("J102").Value = 1733.3334
$ws.Range("K102").Value = 2052.889
$ws.Range("L102").Value = 1733.3334
$ws.Range("M102").Value = -430.8890000000001
$ws.Range("N102").Value = -4977.3334

# Row 132
$ws.Range("H132").Value = 2412.55
$ws.Range("I132").Value = 1950.2
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 5850.6
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -3320.6
$ws.Range("N132").Value = -16458.8

# Row 135
$ws.Range("H135").Value = 34499
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 34499
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 34499
$ws.Range("N135").Value = -44639
$ws.Range("M135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 948
$ws.Range("I93").Value = 790.2857
$ws.Range("K93").Value = 790.2857
$ws.Range("M93").Value = 457.7143

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 39064390
$ws.Range("I122").Value = 44644748
$ws.Range("K122").Value = 133934244
$ws.Range("M122").Value = -133931794
